$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "EISU2205930"
$ws.Range("B3").Value = "EVER ELITE"

# C3 and E3 look like numbers (leading zero / long digit string), so force
# them to be stored as text the same way the other numeric-looking values
# in the sheet (row 2) are stored: as shared strings, not numeric cells.
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "01082"
$ws.Range("C3").ClearFormats()

$ws.Range("D3").Value = "DJLAXA3986096"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "7075381741"
$ws.Range("E3").ClearFormats()

$ws.Range("F3").Value = "EGLV147900108463"
